{"js": "// Append five new paragraphs (\"Line 6\" .. \"Line 10\") after the last\n// existing paragraph in the document body, each fully red (RGB FF0000)\n// \u2014 both the run text and the paragraph mark itself.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\n\nfor (let i = 6; i <= 10; i++) {\n  const newPara = anchor.insertParagraph(`Line ${i}`, \"After\");\n\n  // Color the run text red.\n  newPara.font.color = \"#FF0000\";\n\n  // Color the paragraph mark (pilcrow) red as well, so the paragraph's\n  // own rPr (w:pPr/w:rPr) carries the color too.\n  newPara.getRange(\"End\").font.color = \"#FF0000\";\n\n  anchor = newPara;\n}\n\nawait context.sync();\n", "ps1": "# Append five new paragraphs (\"Line 6\" .. \"Line 10\") at the end of the\n# document, each fully red (RGB FF0000 / wdColorRed) \u2014 both the run text\n# and the paragraph mark itself.\n\n$d = $word.ActiveDocument\n\nfor ($i = 6; $i -le 10; $i++) {\n    $rng = $d.Content\n    $rng.Collapse(0)              # wdCollapseEnd\n    $rng.InsertParagraphAfter()   # create the new paragraph mark\n    $rng.Collapse(0)\n    $rng.Move(1, 1) | Out-Null    # wdCharacter: step into the new paragraph\n    $rng.InsertAfter(\"Line $i\")\n\n    # Color the whole paragraph (run text + the paragraph mark) red.\n    $paraCount = $d.Paragraphs.Count\n    $p = $d.Paragraphs.Item($paraCount)\n    $p.Range.Font.Color = 255     # wdColorRed (RGB FF0000)\n}\n"}
